$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: was "Label" in A6 only -> now "Eigen risico" + extra cells
$ws.Range("A6").Value = "Eigen risico"
$ws.Range("B6").Value = "€ 10043  10044"
$ws.Range("H6").Value = "x"

# --- Row 7: was empty -> now new data row
$ws.Range("B7").Value = 10043
$ws.Range("D7").Value = "rechts"
$ws.Range("E7").Value = "verwijderen"

# --- Row 8: was B8 "asd" -> now becomes a data row under Eigen risico
$ws.Range("B8").Value = 10044
$ws.Range("D8").Value = "links"
$ws.Range("E8").Value = "niet verwijderen"

# --- Row 9: was B9 "asd" -> now header row for a new group
$ws.Range("A9").Value = "n/a"
$ws.Range("B9").Value = "10631 10632 10630          10633"
$ws.Range("H9").Value = "x"

# --- Row 10: was B10 "asd" -> now C10 with data, B10 cleared
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 10631
$ws.Range("D10").Value = "links"
$ws.Range("E10").Value = "verwijderen"

# --- Row 11: was B11 "asd" -> now C11 with data, B11 cleared
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 10632
$ws.Range("D11").Value = "links"
$ws.Range("E11").Value = "verwijderen"

# --- Row 12: was B12 "asd" -> now C12 with data, B12 cleared
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 10630
$ws.Range("D12").Value = "links"
$ws.Range("E12").Value = "verwijderen"

# --- Row 13: was B13 "asd" -> now C13 with data, B13 cleared
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 10633
$ws.Range("D13").Value = "links"
$ws.Range("E13").Value = "verwijderen"

# --- Row 14: was B14 "asd" -> row removed entirely
$ws.Rows.Item(14).ClearContents()

# --- Rows 15-21 stay as-is ("asd" in column B)

# --- Rows 22-38: removed entirely
$ws.Range("B22:B38").ClearContents()
